$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates
$ws.Range("A8").Value = "Volume 30   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/6/2023  Through  11/12/2023"

# Plain numeric value updates
$ws.Range("N15").Value = -75.714285714285
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 130
$ws.Range("J16").Value = 138
$ws.Range("K16").Value = -5.797101449275
$ws.Range("L16").Value = -15.032679738562
$ws.Range("M16").Value = -56.521739130434
$ws.Range("N16").Value = -87.974098057354
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -10
$ws.Range("I17").Value = 259
$ws.Range("J17").Value = 290
$ws.Range("K17").Value = -10.689655172413
$ws.Range("L17").Value = -11.301369863013
$ws.Range("M17").Value = -8.480565371024
$ws.Range("N17").Value = -66.187989556135
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -85.714285714285
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = -86.956521739130
$ws.Range("I18").Value = 122
$ws.Range("J18").Value = 186
$ws.Range("K18").Value = -34.408602150537
$ws.Range("L18").Value = -18.120805369127
$ws.Range("M18").Value = -41.626794258373
$ws.Range("N18").Value = -82.743988684582
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -40
$ws.Range("F19").Value = 30
$ws.Range("H19").Value = -21.052631578947
$ws.Range("I19").Value = 296
$ws.Range("J19").Value = 381
$ws.Range("K19").Value = -22.309711286089
$ws.Range("L19").Value = -8.641975308641
$ws.Range("M19").Value = -2.310231023102
$ws.Range("N19").Value = -8.074534161490
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -30.769230769230
$ws.Range("I20").Value = 95
$ws.Range("J20").Value = 123
$ws.Range("K20").Value = -22.764227642276
$ws.Range("L20").Value = 6.741573033707
$ws.Range("M20").Value = 26.666666666666
$ws.Range("N20").Value = -82.007575757575
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -53.846153846153
$ws.Range("F21").Value = 64
$ws.Range("H21").Value = -37.864077669902
$ws.Range("I21").Value = 923
$ws.Range("J21").Value = 1142
$ws.Range("K21").Value = -19.176882661996
$ws.Range("L21").Value = -10.038986354775
$ws.Range("M21").Value = -23.529411764705
$ws.Range("N21").Value = -73.621034581308
$ws.Range("H22").Value = -100
$ws.Range("M22").Value = -52.173913043478
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 75
$ws.Range("J23").Value = 75
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = -21.052631578947
$ws.Range("M23").Value = 2.739726027397
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 62
$ws.Range("G24").Value = 63
$ws.Range("H24").Value = -1.587301587301
$ws.Range("I24").Value = 746
$ws.Range("J24").Value = 772
$ws.Range("K24").Value = -3.367875647668
$ws.Range("L24").Value = 41.555977229601
$ws.Range("M24").Value = 5.218617771509
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -27.272727272727
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 48.275862068965
$ws.Range("I25").Value = 436
$ws.Range("J25").Value = 353
$ws.Range("K25").Value = 23.512747875354
$ws.Range("L25").Value = 45.333333333333
$ws.Range("M25").Value = -39.695712309820
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -100
$ws.Range("J26").Value = 26
$ws.Range("K26").Value = -3.846153846153
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("L27").Value = -32.608695652173
$ws.Range("L28").Value = -60.526315789473
$ws.Range("N28").Value = -90.853658536585
$ws.Range("L29").Value = -58.620689655172
$ws.Range("N29").Value = -91.891891891891

# Cells changing between text-placeholder and numeric (style + type change)
function Set-AsNumberWithFormat($addr, $donor, $value) {
    $ws.Range($donor).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $value
}
function Set-AsTextWithFormat($addr, $donor, $text) {
    $ws.Range("ZZ1").NumberFormat = "@"
    $ws.Range("ZZ1").Value = $text
    $ws.Range($donor).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range("ZZ1").Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $ws.Range("ZZ1").Clear()
}

Set-AsNumberWithFormat "C16" "I14" 1
Set-AsNumberWithFormat "C23" "I14" 1
Set-AsNumberWithFormat "D23" "I14" 3
Set-AsNumberWithFormat "E23" "K14" -66.666666666666
Set-AsNumberWithFormat "D26" "I14" 1
Set-AsNumberWithFormat "E26" "K14" -100
Set-AsTextWithFormat "F22" "C14" "0"
Set-AsTextWithFormat "F26" "C14" "0"
